$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at A, shifting all existing columns (A:AC) to (B:AD)
$ws.Columns("A:A").Insert()

# Populate the new "Match ID" column
$ws.Range("A3").Value = "Match ID"
$ws.Range("A3").Font.Bold = $true

$ws.Range("A4:A19").Value = 14
$ws.Range("A4:A19").Font.Bold = $true

$ws.Range("A20").Value = 14
$ws.Rows(20).AutoFit()

# Update the active selection to the new Match ID data range
$ws.Range("A3:A19").Select() | Out-Null
